# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows into the Mango dataset:
#   - a new row that becomes row 84 (pushing the former rows 84-187 down by one)
#   - a new row that becomes row 147 (pushing the rows that are now 147-188 down by one more)
# After both inserts the sheet grows from 187 rows (A1:T187) to 189 rows (A1:T189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the blank rows first (this shifts everything below down) ---
$ws.Rows.Item(84).Insert()
$ws.Rows.Item(147).Insert()

# --- Fill in the first new row (final row 84) ---
$ws.Cells.Item(84, 1).Value = 4
$ws.Cells.Item(84, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(84, 3).Value = "Los Lagos"
$ws.Cells.Item(84, 4).Value = 44664
$ws.Cells.Item(84, 5).Value = 10
$ws.Cells.Item(84, 6).Value = "Fruta"
$ws.Cells.Item(84, 7).Value = 100108
$ws.Cells.Item(84, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(84, 9).Value = 100108002
$ws.Cells.Item(84, 10).Value = "Mango"
$ws.Cells.Item(84, 11).Value = "Sin especificar"
$ws.Cells.Item(84, 12).Value = "Primera"
$ws.Cells.Item(84, 13).Value = 30
$ws.Cells.Item(84, 14).Value = 8000
$ws.Cells.Item(84, 15).Value = 8500
$ws.Cells.Item(84, 16).Value = 8250
$ws.Cells.Item(84, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(84, 18).Value = "Perú"
$ws.Cells.Item(84, 19).Value = 2062
$ws.Cells.Item(84, 20).Value = 4

# --- Fill in the second new row (final row 147) ---
$ws.Cells.Item(147, 1).Value = 4
$ws.Cells.Item(147, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(147, 3).Value = "Los Lagos"
$ws.Cells.Item(147, 4).Value = 44663
$ws.Cells.Item(147, 5).Value = 10
$ws.Cells.Item(147, 6).Value = "Fruta"
$ws.Cells.Item(147, 7).Value = 100108
$ws.Cells.Item(147, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(147, 9).Value = 100108002
$ws.Cells.Item(147, 10).Value = "Mango"
$ws.Cells.Item(147, 11).Value = "Sin especificar"
$ws.Cells.Item(147, 12).Value = "Primera"
$ws.Cells.Item(147, 13).Value = 200
$ws.Cells.Item(147, 14).Value = 8000
$ws.Cells.Item(147, 15).Value = 8000
$ws.Cells.Item(147, 16).Value = 8000
$ws.Cells.Item(147, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(147, 18).Value = "Perú"
$ws.Cells.Item(147, 19).Value = 2000
$ws.Cells.Item(147, 20).Value = 4
